$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.041937413968151
$ws.Range("D2").Value = 1.048921179698696
$ws.Range("E2").Value = 1.04963823194174
$ws.Range("F2").Value = 1.059264189852732
$ws.Range("I2").Value = 1.038913383302503
$ws.Range("J2").Value = 1.047015931704176
$ws.Range("K2").Value = 1.051679677902946
$ws.Range("L2").Value = 1.05239473330045
$ws.Range("M2").Value = 1.061994168364726

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.042962850435983
$ws.Range("D3").Value = 1.049729714893526
$ws.Range("E3").Value = 1.050538416031361
$ws.Range("F3").Value = 1.060224051541916
$ws.Range("I3").Value = 1.039121596199218
$ws.Range("J3").Value = 1.047687223520756
$ws.Range("K3").Value = 1.0523003648002
$ws.Range("L3").Value = 1.053106976685382
$ws.Range("M3").Value = 1.062767856678568

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.043626590776393
$ws.Range("D4").Value = 1.050252670218264
$ws.Range("E4").Value = 1.051121406435495
$ws.Range("F4").Value = 1.060845543388529
$ws.Range("I4").Value = 1.039254491063155
$ws.Range("J4").Value = 1.048121215516117
$ws.Range("K4").Value = 1.052701117784367
$ws.Range("L4").Value = 1.053567722648926
$ws.Range("M4").Value = 1.063268258465172

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.043905678106922
$ws.Range("D5").Value = 1.050472466664375
$ws.Range("E5").Value = 1.051366617006588
$ws.Range("F5").Value = 1.061106912647685
$ws.Range("I5").Value = 1.039309920723711
$ws.Range("J5").Value = 1.048303574409506
$ws.Range("K5").Value = 1.052869384522391
$ws.Range("L5").Value = 1.053761390012927
$ws.Range("M5").Value = 1.063478572337808

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.043952541069545
$ws.Range("D6").Value = 1.050509368258875
$ws.Range("E6").Value = 1.051407796031906
$ws.Range("F6").Value = 1.061150803184775
$ws.Range("I6").Value = 1.039319201828029
$ws.Range("J6").Value = 1.048334187893977
$ws.Range("K6").Value = 1.052897624920074
$ws.Range("L6").Value = 1.053793905820107
$ws.Range("M6").Value = 1.06351388170975

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.043630319755534
$ws.Range("D7").Value = 1.050255607362556
$ws.Range("E7").Value = 1.051124682475314
$ws.Range("F7").Value = 1.060849035448919
$ws.Range("I7").Value = 1.039255233444016
$ws.Range("J7").Value = 1.048123652564347
$ws.Range("K7").Value = 1.05270336699825
$ws.Range("L7").Value = 1.053570310561504
$ws.Range("M7").Value = 1.063271068906791

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.042283920376825
$ws.Range("D8").Value = 1.049194473078919
$ws.Range("E8").Value = 1.049942346969431
$ws.Range("F8").Value = 1.059588496777506
$ws.Range("I8").Value = 1.038984128982471
$ws.Range("J8").Value = 1.047242875920609
$ws.Range("K8").Value = 1.051889622195848
$ws.Range("L8").Value = 1.052635464326353
$ws.Range("M8").Value = 1.062255686288692

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.039913053056871
$ws.Range("D9").Value = 1.047322972569441
$ws.Range("E9").Value = 1.047862876926247
$ws.Range("F9").Value = 1.057370351588349
$ws.Range("I9").Value = 1.038492395725084
$ws.Range("J9").Value = 1.045687964197981
$ws.Range("K9").Value = 1.050449050446938
$ws.Range("L9").Value = 1.050987235967661
$ws.Range("M9").Value = 1.060464757017261

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.038333613503127
$ws.Range("D10").Value = 1.046074259278775
$ws.Range("E10").Value = 1.046479277934817
$ws.Range("F10").Value = 1.055893718665013
$ws.Range("I10").Value = 1.038155184264178
$ws.Range("J10").Value = 1.044649464457852
$ws.Range("K10").Value = 1.049484245685973
$ws.Range("L10").Value = 1.049887848409166
$ws.Range("M10").Value = 1.059269713294593

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.037649972640794
$ws.Range("D11").Value = 1.04553331665542
$ws.Range("E11").Value = 1.045880819022422
$ws.Range("F11").Value = 1.055254836878152
$ws.Range("I11").Value = 1.038006947340259
$ws.Range("J11").Value = 1.044199341358543
$ws.Range("K11").Value = 1.049065434719821
$ws.Range("L11").Value = 1.049411675128484
$ws.Range("M11").Value = 1.058751997162736

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.037396078287037
$ws.Range("D12").Value = 1.045332351104643
$ws.Range("E12").Value = 1.045658622996086
$ws.Range("F12").Value = 1.055017605194287
$ws.Range("I12").Value = 1.037951552149934
$ws.Range("J12").Value = 1.044032079146789
$ws.Range("K12").Value = 1.048909713490599
$ws.Range("L12").Value = 1.049234784081913
$ws.Range("M12").Value = 1.058559656845395

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037450537660762
$ws.Range("D13").Value = 1.045375460483694
$ws.Range("E13").Value = 1.045706280338424
$ws.Range("F13").Value = 1.05506848867729
$ws.Range("I13").Value = 1.037963449685944
$ws.Range("J13").Value = 1.044067960457493
$ws.Range("K13").Value = 1.048943123279866
$ws.Range("L13").Value = 1.049272728665671
$ws.Range("M13").Value = 1.058600916181685

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.037628984830982
$ws.Range("D14").Value = 1.045516705496688
$ws.Range("E14").Value = 1.045862450220375
$ws.Range("F14").Value = 1.055235225648229
$ws.Range("I14").Value = 1.03800237515512
$ws.Range("J14").Value = 1.044185516764441
$ws.Range("K14").Value = 1.049052565936719
$ws.Range("L14").Value = 1.04939705364519
$ws.Range("M14").Value = 1.058736099012873

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.037738937368544
$ws.Range("D15").Value = 1.045603726550512
$ws.Range("E15").Value = 1.04595868468373
$ws.Range("F15").Value = 1.05533796808786
$ws.Range("I15").Value = 1.03802631425663
$ws.Range("J15").Value = 1.044257938295536
$ws.Range("K15").Value = 1.049119976493422
$ws.Range("L15").Value = 1.049473651857536
$ws.Range("M15").Value = 1.058819384669622

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.038378990136134
$ws.Range("D16").Value = 1.046110154848391
$ws.Range("E16").Value = 1.046519009416635
$ws.Range("F16").Value = 1.0559361299618
$ws.Range("I16").Value = 1.038164975481625
$ws.Range("J16").Value = 1.044679328297454
$ws.Range("K16").Value = 1.049512018849941
$ws.Range("L16").Value = 1.049919447750462
$ws.Range("M16").Value = 1.059304067112072

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.038780549925383
$ws.Range("D17").Value = 1.046427759921836
$ws.Range("E17").Value = 1.046870660540454
$ws.Range("F17").Value = 1.05631147823753
$ws.Range("I17").Value = 1.038251359397769
$ws.Range("J17").Value = 1.044943536121918
$ws.Range("K17").Value = 1.049757657599731
$ws.Range("L17").Value = 1.050199049049074
$ws.Range("M17").Value = 1.059608028020906

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.03901479869363
$ws.Range("D18").Value = 1.04661299022295
$ws.Range("E18").Value = 1.047075835215555
$ws.Range("F18").Value = 1.056530461689556
$ws.Range("I18").Value = 1.038301531240401
$ws.Range("J18").Value = 1.045097601021328
$ws.Range("K18").Value = 1.049900833737051
$ws.Range("L18").Value = 1.050362123054215
$ws.Range("M18").Value = 1.059785298691234

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.039094675829281
$ws.Range("D19").Value = 1.046676144977944
$ws.Range("E19").Value = 1.047145805057792
$ws.Range("F19").Value = 1.056605137699445
$ws.Range("I19").Value = 1.03831860215655
$ws.Range("J19").Value = 1.045150125851532
$ws.Range("K19").Value = 1.049949635983421
$ws.Range("L19").Value = 1.050417724902246
$ws.Range("M19").Value = 1.059845739216901

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.038737463696626
$ws.Range("D20").Value = 1.046393686313297
$ws.Range("E20").Value = 1.046832925229721
$ws.Range("F20").Value = 1.056271201829547
$ws.Range("I20").Value = 1.038242113397695
$ws.Range("J20").Value = 1.04491519357219
$ws.Range("K20").Value = 1.049731313306717
$ws.Range("L20").Value = 1.05016905178508
$ws.Range("M20").Value = 1.059575418429773

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.037576435472864
$ws.Range("D21").Value = 1.04547511333003
$ws.Range("E21").Value = 1.045816459363561
$ws.Range("F21").Value = 1.055186123621198
$ws.Range("I21").Value = 1.037990921774387
$ws.Range("J21").Value = 1.044150901192154
$ws.Range("K21").Value = 1.049020342108194
$ws.Range("L21").Value = 1.049360443559376
$ws.Range("M21").Value = 1.058696292067761

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036846684221822
$ws.Range("D22").Value = 1.04489736581571
$ws.Range("E22").Value = 1.04517793583522
$ws.Range("F22").Value = 1.054504340443097
$ws.Range("I22").Value = 1.037831058923695
$ws.Range("J22").Value = 1.043669976035891
$ws.Range("K22").Value = 1.04857242297051
$ws.Range("L22").Value = 1.048851929341111
$ws.Range("M22").Value = 1.058143333849168

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.037233516987274
$ws.Range("D23").Value = 1.045203659759371
$ws.Range("E23").Value = 1.045516375041026
$ws.Range("F23").Value = 1.054865723809981
$ws.Range("I23").Value = 1.037915987914052
$ws.Range("J23").Value = 1.043924959800976
$ws.Range("K23").Value = 1.048809958779165
$ws.Range("L23").Value = 1.049121512581223
$ws.Range("M23").Value = 1.058436487755638

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.038756932429144
$ws.Range("D24").Value = 1.04640908278253
$ws.Range("E24").Value = 1.046849975999171
$ws.Range("F24").Value = 1.056289400852768
$ws.Range("I24").Value = 1.038246291929704
$ws.Range("J24").Value = 1.044928000483431
$ws.Range("K24").Value = 1.049743217470171
$ws.Range("L24").Value = 1.050182606296702
$ws.Range("M24").Value = 1.059590153376303

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.040525780019804
$ws.Range("D25").Value = 1.04780698913746
$ws.Range("E25").Value = 1.048399995321741
$ws.Range("F25").Value = 1.057943423874374
$ws.Range("I25").Value = 1.038621177808478
$ws.Range("J25").Value = 1.046090282396493
$ws.Range("K25").Value = 1.050822255890437
$ws.Range("L25").Value = 1.051413445874499
$ws.Range("M25").Value = 1.060927951005201

